# "add data until november 15th"
#
# Appends daily COVID bulletin rows 495-502 (2021-11-08 .. 2021-11-15) to
# "Planilha1", following the sheet's existing layout:
#   - columns A-J are literal values straight from the bulletin
#   - columns K-P are "value minus same column in the previous row" deltas
# The row directly above each new row is duplicated first (Copy +
# PasteSpecial formats-only) so the new cells pick up the same style
# indices already used throughout the sheet instead of creating new ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, A(date serial), B, C, D, E, F, G, H, I, J
$newRows = @(
    @(495, 44508, 17119, 21, 7048, 24188, 6889, 5,  3, 2,  154),
    @(496, 44509, 17152, 26, 7048, 24226, 6889, 5,  3, 2,  154),
    @(497, 44510, 17171, 23, 7050, 24244, 6889, 7,  1, 6,  154),
    @(498, 44511, 17198, 30, 7053, 24281, 6889, 10, 1, 9,  154),
    @(499, 44512, 17227, 24, 7054, 24305, 6889, 11, 1, 10, 154),
    @(500, 44513, 17240, 27, 7057, 24324, 6889, 14, 1, 13, 154),
    @(501, 44514, 17240, 27, 7057, 24324, 6889, 14, 1, 13, 154),
    @(502, 44515, 17249, 15, 7060, 24324, 6891, 15, 1, 14, 154)
)

foreach ($entry in $newRows) {
    $r = $entry[0]
    $prev = $r - 1

    # Clone formatting from the row above so we reuse existing style ids.
    $ws.Range("A$prev`:P$prev").Copy() | Out-Null
    $ws.Range("A$r`:P$r").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
    $ws.Cells.Item($r, 5).Value = $entry[5]
    $ws.Cells.Item($r, 6).Value = $entry[6]
    $ws.Cells.Item($r, 7).Value = $entry[7]
    $ws.Cells.Item($r, 8).Value = $entry[8]
    $ws.Cells.Item($r, 9).Value = $entry[9]
    $ws.Cells.Item($r, 10).Value = $entry[10]

    $ws.Cells.Item($r, 11).Formula = "=D$r-D$prev"
    $ws.Cells.Item($r, 12).Formula = "=F$r-F$prev"
    $ws.Cells.Item($r, 13).Formula = "=B$r-B$prev"
    $ws.Cells.Item($r, 14).Formula = "=J$r-J$prev"
    $ws.Cells.Item($r, 15).Formula = "=G$r-G$prev"
    $ws.Cells.Item($r, 16).Formula = "=C$r-C$prev"
}

$excel.CutCopyMode = $false

# Move the active cell/selection to match where the author left off editing.
$ws.Range("M510").Select()
